$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; existing D:K data shifts to F:M
$ws.Columns("D:E").Insert()

# New D/E columns should carry the same number formatting as the (shifted) F/G columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (latest two quarters) with their reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 84400
$ws.Range("E8").Value = 78500
$ws.Range("D9").Value = 13100
$ws.Range("E9").Value = 8700
$ws.Range("D10").Value = 71300
$ws.Range("E10").Value = 69800
$ws.Range("D12").Value = 146200
$ws.Range("E12").Value = 86600
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 200
$ws.Range("E15").Value = 200
$ws.Range("D17").Value = 223800
$ws.Range("E17").Value = 148600
$ws.Range("D18").Value = -139400
$ws.Range("E18").Value = -70100
$ws.Range("D20").Value = -2300
$ws.Range("E20").Value = -7000
$ws.Range("D21").Value = -138100
$ws.Range("E21").Value = -73500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -141700
$ws.Range("E23").Value = -77100
$ws.Range("D24").Value = -800
$ws.Range("E24").Value = -700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -140900
$ws.Range("E26").Value = -76400
$ws.Range("D27").Value = -140900
$ws.Range("E27").Value = -76400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2300
$ws.Range("E32").Value = 7000
$ws.Range("D33").Value = -140900
$ws.Range("E33").Value = -76400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -140900
$ws.Range("E35").Value = -76400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 370800
$ws.Range("E41").Value = 209700
$ws.Range("D42").Value = 803100
$ws.Range("E42").Value = 583200
$ws.Range("D43").Value = 66600
$ws.Range("E43").Value = 59400
$ws.Range("D44").Value = 125400
$ws.Range("E44").Value = 115800
$ws.Range("D45").Value = 60300
$ws.Range("E45").Value = 44000
$ws.Range("D46").Value = 1426200
$ws.Range("E46").Value = 1012100
$ws.Range("D47").Value = 31700
$ws.Range("E47").Value = 30000
$ws.Range("D48").Value = 97000
$ws.Range("E48").Value = 76800
$ws.Range("D49").Value = 11600
$ws.Range("E49").Value = 15300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 75600
$ws.Range("E52").Value = 48700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1642100
$ws.Range("E54").Value = 1182900
$ws.Range("D57").Value = 33800
$ws.Range("E57").Value = 20400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 139900
$ws.Range("E59").Value = 94000
$ws.Range("D60").Value = 173700
$ws.Range("E60").Value = 114400
$ws.Range("D61").Value = 420600
$ws.Range("E61").Value = 415400
$ws.Range("D62").Value = 15600
$ws.Range("E62").Value = 13200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 609800
$ws.Range("E66").Value = 543100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1578900
$ws.Range("E72").Value = -1438000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1032300
$ws.Range("E76").Value = 639800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -140900
$ws.Range("E81").Value = -76400
$ws.Range("D83").Value = 3500
$ws.Range("E83").Value = 3600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -122400
$ws.Range("E89").Value = -113300
$ws.Range("D91").Value = -20200
$ws.Range("E91").Value = -20100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -237700
$ws.Range("E94").Value = -63400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 521300
$ws.Range("E100").Value = -24000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 161100
$ws.Range("E102").Value = -200700

# Row 91 (Capital Expenditures) was restated for several historical quarters
$ws.Range("F91").Value = -8700
$ws.Range("G91").Value = -12200
$ws.Range("H91").Value = -3900
$ws.Range("I91").Value = -800
$ws.Range("J91").Value = -2900
